# Fill in the quantity/count column (G) for the exam-duty bill rows.
# These are the only user-entered inputs; the dependent formulas in
# column I (and the grand-total in I32) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G9").Value = 117
$ws.Range("G12").Value = 117
$ws.Range("G14").Value = 119
$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
$ws.Range("G29").Value = 15
